$d = $word.ActiveDocument

# 1. "a architecture" -> "an architecture"
$d.Content.Find.Execute("a architecture", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "an architecture", 2)

# 2. "The  markup" (double space) -> "The markup" (single space)
$d.Content.Find.Execute("The  markup", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "The markup", 2)

# 3. "can be persist in" -> "can be persisted in"
$d.Content.Find.Execute("can be persist in", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "can be persisted in", 2)

# 4. "has to modified" -> "has to be modified"
$d.Content.Find.Execute("has to modified", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "has to be modified", 2)

Write-Output "done"
